$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# J14 holds a standalone formula (ratio of G2/E2)
$ws.Range("J14").Formula = "=G2/E2"

# J15:J23 were originally filled as one shared formula (G3/E3 relative pattern),
# dragged down to row 23, so the shared formula's range spans J15:J23 even
# though rows 22:23 end up cleared afterwards.
$ws.Range("J15:J23").Formula = "=G3/E3"
$ws.Range("J22:J23").ClearContents()

# Restore the active selection to E8
$ws.Range("E8").Select()
